$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only these columns vary per-row in this dataset (the rest -A,B,C,E,F,G,H,N,Q,R-
# are identical for every data row): D=Fecha, I=Calidad, J=Volumen,
# K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado, O=Origen,
# P=Precio $/Kg.
$cols = @(4, 9, 10, 11, 12, 13, 15, 16)

# The data rows (2..43) were reshuffled: each destination row now carries the
# values that used to live in another (source) row. Row 32 is unchanged.
$map = @{
    2 = 42; 3 = 36; 4 = 2;  5 = 19; 6 = 37; 7 = 28; 8 = 10; 9 = 11; 10 = 4;
    11 = 31; 12 = 25; 13 = 33; 14 = 34; 15 = 20; 16 = 3; 17 = 30; 18 = 5; 19 = 23;
    20 = 43; 21 = 26; 22 = 9; 23 = 6; 24 = 13; 25 = 29; 26 = 39; 27 = 17; 28 = 24;
    29 = 22; 30 = 38; 31 = 18; 32 = 32; 33 = 40; 34 = 12; 35 = 15; 36 = 16; 37 = 21;
    38 = 7; 39 = 35; 40 = 27; 41 = 14; 42 = 41; 43 = 8
}

# Snapshot every original value first so the permutation can be applied
# without a source row being overwritten before it has been read.
$snapshot = @{}
for ($r = 2; $r -le 43; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
